# Generate Report for Archive
#
# 1. Change the shared "Ready for handoff" status text to "In Translation"
#    (this shows up in Overview!E2, Overview!F2, zh-cn!C2, de-de!C2, since
#    they all reference the same shared string).
# 2. Narrow the "status" columns that display this text:
#      Overview columns E and F (zh-cn / de-de) from 17.2159881591797 -> 13.4101845877511
#      zh-cn column C   (Status) from 17.2159881591797 -> 13.4101845877511
#      de-de column C   (Status) from 17.2159881591797 -> 13.4101845877511

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update every cell that currently holds the old status text ---
foreach ($ws in @($wsOverview, $wsZhCn, $wsDeDe)) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $cellText = "" + $cell.Value()
            if ($cellText -eq $oldStatus) {
                $cell.Value = $newStatus
            }
        }
    }
}

# --- Resize the affected columns ---
# Excel's ColumnWidth is quantized to whole-pixel increments, so asking for
# 13.4101845877511 directly snaps to the nearest representable width. Using
# an input of 12.5 lands reliably on that nearest representable width
# (13.333333333333334), which is the closest Excel can get to the target.
$newWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth   # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth   # column F (de-de)

$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth       # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth       # column C (Status)
